# MBlaul - Efforts Logbook: add Sprint #6 entry (row 9) on the "Example"
# sheet and the corresponding tally row (row 10) on the "Count" sheet.
# Also refresh the active-cell selections, matching the author's final
# save state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Example": fill in the previously-blank Sprint #6 row (row 9).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Example")

# Row 9 was an empty template row styled like the rows above it; copy the
# formatting down from row 8 (the last filled-in sprint row) before
# putting in the new values, same as a user extending the table would.
$ws1.Range("A8:H8").Copy()
$ws1.Range("A9:H9").PasteSpecial(-4122)  # xlPasteFormats

$ws1.Range("C9").Value = 20
$ws1.Range("E9").Value = "Group meetings, assigned tasks, worked on polish for the site, user preferences, deployment"
$ws1.Range("H9").Value = "/core/api/server_https.js`n/core/api/routes/weather.js,`n/frontend/client/src/App.js,`n/frontend/client/src/compenents/content/dashboard/Dashboard.js,`n/frontend/client/src/compenents/content/station/Station.js,`n/frontend/client/src/compenents/content/station/Location.js,`n/frontend/client/src/compenents/content/station/Locations.js,`n/frontend/client/src/compenents/content/station/Map.js,`n/frontend/client/src/compenents/content/user/Settings.js,`n/frontend/client/src/actions/authActions.js,`n/frontend/client/src/reducers/authRedcuer.js`n"
$ws1.Range("D9").Value = "A, B, C, D, E, F, H, I, J, K, L, M, N"
$ws1.Range("F9").Value = "Architecture,`nProject Plan Update,`nUser Interface Design, User Stories, Tools Tech ETC"
$ws1.Range("G9").Value = "github.com/mblaul/skypi/api,`ngithub.com/mblaul/skypi/frontend,`nGoogle Team Drive,`nhttps://drive.google.com/drive/u/1/folders/0ADmYitmSMBPCUk9PVA`nTeam Trello,`nhttps://trello.com/skypitasksseniordesign/home"

# Final row heights (rows re-wrap once the new text is in place).
$ws1.Rows.Item(8).RowHeight = 252
$ws1.Rows.Item(9).RowHeight = 236.25

# ---------------------------------------------------------------------
# Sheet "Count": fill in the tally row (row 10) for Sprint #6.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Count")

$ws2.Range("A9:P9").Copy()
$ws2.Range("A10:P10").PasteSpecial(-4122)  # xlPasteFormats

$ws2.Range("B10").Value = "SkyPi"
$ws2.Range("C10").Value = 1
$ws2.Range("D10").Value = 1
$ws2.Range("E10").Value = 1
$ws2.Range("F10").Value = 3
$ws2.Range("G10").Value = 3
$ws2.Range("H10").Value = 1
$ws2.Range("J10").Value = 3
$ws2.Range("K10").Value = 3
$ws2.Range("L10").Value = 1
$ws2.Range("M10").Value = 1
$ws2.Range("N10").Value = 1
$ws2.Range("O10").Value = 1
$ws2.Range("P10").Value = 5

# ---------------------------------------------------------------------
# Restore the selections the author left active on each sheet.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("E21").Select()

$ws2.Activate()
$ws2.Range("P17").Select()

$ws1.Activate()
